$wb = $excel.ActiveWorkbook

# --- Sheet 1: LH-TC-ADMINCONSTRAINS-Review ---
# Fill in "Reviewer verification" (column J) as Closed for rows 2-4
$wsReview = $wb.Worksheets.Item("LH-TC-ADMINCONSTRAINS-Review")
$wsReview.Range("J2").Value = "Closed"
$wsReview.Range("J3").Value = "Closed"
$wsReview.Range("J4").Value = "Closed"

# --- Sheet 2: Version History ---
# Add a new version history row (v1.2) documenting that the comments were closed
$wsHistory = $wb.Worksheets.Item("Version History")
$wsHistory.Range("A4").Value = "v1.2"
$wsHistory.Range("B4").Value = "Omar Sherif"
$wsHistory.Range("C4").Value = "The comments closed "
$wsHistory.Range("D4").Formula = "=DATE(2025,4,22)"
